# Applies the "additional scraping" commit:
#   1. Inserts a new "Player Info" sheet at the front of the workbook with
#      ID / NAME / BATTING_HAND / BOWL_STYLE columns for player 4592.
#   2. Renames the MATCH_CARD_LINK column to MATCH_CODE on both the
#      "ODI Batting" and "ODI Bowling" sheets.
#   3. Replaces the full howstat.com scorecard URL values in that column
#      with just the numeric match code that was at the end of the URL.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, placed before "ODI Batting".
#    Copying an existing sheet (rather than Worksheets.Add()) keeps the
#    same sheetPr / sheetFormatPr / pageMargins conventions used by the
#    rest of the workbook, then we wipe it and fill in the player data.
#
#    NOTE: worksheet handles in this host are position-based, not
#    identity-based, so `Copy()` (which shifts everybody's index) quietly
#    repoints any variable captured beforehand. Sheets referenced after a
#    structural change (copy/add/move/delete) are always re-fetched by
#    name below instead of reusing older handles.
# ---------------------------------------------------------------------
$battingForCopy = $wb.Worksheets.Item("ODI Batting")
$battingForCopy.Copy($battingForCopy)
$playerInfo = $wb.Worksheets.Item(1)
$playerInfo.Name = "Player Info"
$playerInfo.Cells.Clear()

# Reuse the workbook's existing bold/bordered/centred header style (rather
# than rebuilding it with Font/Borders/Alignment, which would fork a new,
# near-duplicate style resource) by copy/pasting the format from one of
# the other sheets' header cells.
$headerStyleSource = $wb.Worksheets.Item("ODI Batting")
$headerStyleSource.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Player id looks numeric ("4592") but the source data stores it (like
# every other column in this workbook) as text, so force a text format
# before writing it in.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4592"
$playerInfo.Range("B2").Value = "Kesrick Omari Kenal Williams"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

[void]$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------
# 2 & 3. MATCH_CARD_LINK -> MATCH_CODE, url -> trailing numeric code.
#    Re-fetch these two sheets by name now that the sheet collection has
#    been reshuffled by the copy above.
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")
function Convert-LinkColumnToCode($sheet, $headerCell, $firstDataRow, $lastDataRow, $col) {
    $sheet.Range($headerCell).Value = "MATCH_CODE"

    for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
        $cell = $sheet.Cells.Item($r, $col)
        $link = $cell.Value()
        if ($link) {
            $code = $link.Substring($link.LastIndexOf("=") + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }
}

# ODI Batting: MATCH_CARD_LINK is column D (4), data rows 2-9.
Convert-LinkColumnToCode $batting "D1" 2 9 4 | Out-Null

# ODI Bowling: MATCH_CARD_LINK is column B (2), data rows 2-9.
Convert-LinkColumnToCode $bowling "B1" 2 9 2 | Out-Null
